$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test-case rows appended below the existing "list price too low" check,
# entered in the order the author actually typed them.

$ws.Range("A12").Value = "102_AutomobileInsurance_002_VehicleData_002_EnterNumericValuesBelowRange"
$ws.Range("A12").Style = "Standard"

$ws.Range("H12").Value = "499.9999"

$ws.Range("A11").Value = "Vehicle Page check error hint list value ranges"

$ws.Range("D12").Value = "0"

$ws.Range("D11").Value = "<HINT Must be a number between 1 and 2000>"

$ws.Range("J11").Value = "<HINT Must be a number between 100 and 100000>"

$ws.Range("J12").Value = "-8888"

$ws.Range("A13").Value = "102_AutomobileInsurance_002_VehicleData_002_EnterNumericValuesAboveRange"
$ws.Range("A13").Style = "Standard"

$ws.Range("D13").Value = "2001"

$ws.Range("H13").Value = "999999999999999"

$ws.Range("J13").Value = "100000.0001"

$ws.Range("E15").Value = "<HINT Must be today or somewhere in the past>"

$ws.Range("A14").Value = "102_AutomobileInsurance_002_VehicleData_002_ManufacturingDateInTheFuture"
$ws.Range("A14").Style = "Standard"

$ws.Range("A15").Value = "Vehicle Page check error hint manufacturing date in the future"
$ws.Range("A15").Style = "Standard"

$ws.Range("E14").Value = "12/31/2099"

# H11 keeps the pre-existing "500 and 100000" hint text (unchanged value, but it is
# the only still-referenced string from the old row 11, so Excel keeps it in place).
$ws.Range("H11").Value = "<HINT Must be a number between 500 and 100000>"

# Move the selection to the last-edited cell, as Excel does after typing.
$ws.Range("A15").Select()

# The new D/E/J columns needed re-measuring (bestFit) once the new, longer hint
# strings were entered, producing the final autofit column widths.
$ws.Columns("D").ColumnWidth = 39.6640625
$ws.Columns("E").ColumnWidth = 40.88671875
$ws.Columns("J").ColumnWidth = 43.77734375

# The picture (anchored one-cell, not resizing with cells) is pushed down by the
# four newly-inserted rows above it.
$pic = $ws.Shapes.Item(1)
$pic.Left = 38100 / 914400 * 72
$pic.Top = 2926080 / 914400 * 72
